$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" (Changed) date column C for rows 2-18 moves forward by
# one day: serial date 45205 (2023-10-06) -> 45206 (2023-10-07).
foreach ($row in 2..18) {
    $ws.Cells.Item($row, 3).Value = 45206
}
